$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Peak Demand (kWh) values for days 1-19 (rows 3-21)
$values = @(
    1.5972290185,
    1.4333769443,
    1.4648810306,
    1.3283572192,
    1.5178888889,
    1.3748507687,
    1.3417896769,
    1.6333748958,
    1.4929037057,
    1.5098210489,
    1.4353980555,
    1.3989891665,
    1.5824213889,
    1.3833822222,
    1.3633630558,
    1.6230982788,
    1.365,
    1.5436086043,
    1.5989197218
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Add new row 22 for day 20
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = 1.4707985559
